$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rename the "15to30" product-table id fragment to "15to45" everywhere it
#    occurs. The fragment is used inside several Jinja-like placeholders,
#    e.g. "{{ tblProducts_1_15to30.labels.kv.idx }}" and
#    "{%tr for prod in tblProducts_1_15to30.tb_items %}".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("15to30", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "15to45", 2)

# ---------------------------------------------------------------------------
# 2. Tweak the two "discount" rows of the products table: make the row a bit
#    shorter and center a few of the cells that used to be right aligned (or
#    unaligned).
# ---------------------------------------------------------------------------

# Locate the products table (7 columns, whose first cell references
# "tblProducts") instead of hard-coding its position among $d.Tables.
$productsTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(1, 1).Range.Text -like "*tblProducts*") {
        $productsTable = $candidate
        break
    }
}

foreach ($rowIdx in 3, 4) {
    $r = $productsTable.Rows.Item($rowIdx)

    # Row height: 432 -> 288 twips (21.6pt -> 14.4pt).
    $r.Height = 14.4

    # wdAlignParagraphLeft = 0, wdAlignParagraphCenter = 1

    # Cell 1 (idx column): right -> center.
    $r.Cells.Item(1).Range.Paragraphs.Item(1).Format.Alignment = 1
    # Cell 2 (sku column): unaligned/left -> center.
    $r.Cells.Item(2).Range.Paragraphs.Item(1).Format.Alignment = 1
    # Cell 3 (name column): center -> unaligned/left (the default).
    $r.Cells.Item(3).Range.Paragraphs.Item(1).Format.Alignment = 0
    # Cell 4 (qty column): right -> center.
    $r.Cells.Item(4).Range.Paragraphs.Item(1).Format.Alignment = 1
    # Cell 5 (tax/discount amount column): right -> center.
    $r.Cells.Item(5).Range.Paragraphs.Item(1).Format.Alignment = 1
}
